$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) and C (Link) updates ---
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("B39").Value = 'Frax'
$ws.Range("C39").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

# --- Column D (Price) updates, staged through a helper column (Z) ---
# Using PasteSpecial(xlPasteValues) keeps these as plain text without
# Excel auto-converting numeric-looking strings (e.g. "279.93") into
# numbers, and without altering the destination cell style.
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = '20.836.44'
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = '1.484.03'
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = '1.007'
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = '1.006'
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = '279.93'
$ws.Range("Z7").NumberFormat = "@"
$ws.Range("Z7").Value = '0.3780'
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = '0.3089'
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = '42.21'
$ws.Range("Z10").NumberFormat = "@"
$ws.Range("Z10").Value = '0.06759'
$ws.Range("Z11").NumberFormat = "@"
$ws.Range("Z11").Value = '1.022'
$ws.Range("Z12").NumberFormat = "@"
$ws.Range("Z12").Value = '1.007'
$ws.Range("Z13").NumberFormat = "@"
$ws.Range("Z13").Value = '5.478'
$ws.Range("Z14").NumberFormat = "@"
$ws.Range("Z14").Value = '17.68'
$ws.Range("Z15").NumberFormat = "@"
$ws.Range("Z15").Value = '1.491.76'
$ws.Range("Z16").NumberFormat = "@"
$ws.Range("Z16").Value = '6.268'
$ws.Range("Z17").NumberFormat = "@"
$ws.Range("Z17").Value = '0.00001040'
$ws.Range("Z18").NumberFormat = "@"
$ws.Range("Z18").Value = '0.06494'
$ws.Range("Z19").NumberFormat = "@"
$ws.Range("Z19").Value = '80.23'
$ws.Range("Z20").NumberFormat = "@"
$ws.Range("Z20").Value = '1.007'
$ws.Range("Z21").NumberFormat = "@"
$ws.Range("Z21").Value = '5.881'
$ws.Range("Z22").NumberFormat = "@"
$ws.Range("Z22").Value = '14.84'
$ws.Range("Z23").NumberFormat = "@"
$ws.Range("Z23").Value = '10.85'
$ws.Range("Z24").NumberFormat = "@"
$ws.Range("Z24").Value = '2.329'
$ws.Range("Z25").NumberFormat = "@"
$ws.Range("Z25").Value = '20.830.32'
$ws.Range("Z26").NumberFormat = "@"
$ws.Range("Z26").Value = '2.291'
$ws.Range("Z27").NumberFormat = "@"
$ws.Range("Z27").Value = '145.37'
$ws.Range("Z28").NumberFormat = "@"
$ws.Range("Z28").Value = '17.67'
$ws.Range("Z29").NumberFormat = "@"
$ws.Range("Z29").Value = '4.796'
$ws.Range("Z30").NumberFormat = "@"
$ws.Range("Z30").Value = '1.657.71'
$ws.Range("Z31").NumberFormat = "@"
$ws.Range("Z31").Value = '112.28'
$ws.Range("Z32").NumberFormat = "@"
$ws.Range("Z32").Value = '5.745'
$ws.Range("Z33").NumberFormat = "@"
$ws.Range("Z33").Value = '0.9261'
$ws.Range("Z34").NumberFormat = "@"
$ws.Range("Z34").Value = '0.07868'
$ws.Range("Z35").NumberFormat = "@"
$ws.Range("Z35").Value = '8.364'
$ws.Range("Z36").NumberFormat = "@"
$ws.Range("Z36").Value = '1.454'
$ws.Range("Z37").NumberFormat = "@"
$ws.Range("Z37").Value = '11.16'
$ws.Range("Z38").NumberFormat = "@"
$ws.Range("Z38").Value = '4.882'
$ws.Range("Z39").NumberFormat = "@"
$ws.Range("Z39").Value = '1.007'
$ws.Range("Z40").NumberFormat = "@"
$ws.Range("Z40").Value = '0.05734'
$ws.Range("Z41").NumberFormat = "@"
$ws.Range("Z41").Value = '0.1948'
$ws.Range("Z42").NumberFormat = "@"
$ws.Range("Z42").Value = '0.02076'
$ws.Range("Z43").NumberFormat = "@"
$ws.Range("Z43").Value = '1.138'
$ws.Range("Z44").NumberFormat = "@"
$ws.Range("Z44").Value = '0.5483'
$ws.Range("Z45").NumberFormat = "@"
$ws.Range("Z45").Value = '3.625'
$ws.Range("Z46").NumberFormat = "@"
$ws.Range("Z46").Value = '12.54'
$ws.Range("Z47").NumberFormat = "@"
$ws.Range("Z47").Value = '0.5299'
$ws.Range("Z48").NumberFormat = "@"
$ws.Range("Z48").Value = '1.819'
$ws.Range("Z49").NumberFormat = "@"
$ws.Range("Z49").Value = '1.106'
$ws.Range("Z50").NumberFormat = "@"
$ws.Range("Z50").Value = '111.52'
$ws.Range("Z51").NumberFormat = "@"
$ws.Range("Z51").Value = '0.06443'
$ws.Range("Z2:Z51").Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$ws.Range("Z2:Z51").Clear()

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = '  -5.74%  '
$ws.Range("E3").Value = '  -4.73%  '
$ws.Range("E4").Value = '  +0.94%  '
$ws.Range("E6").Value = '  -3.76%  '
$ws.Range("E7").Value = '  -5.11%  '
$ws.Range("E8").Value = '  -4.54%  '
$ws.Range("E9").Value = '  -4.06%  '
$ws.Range("E10").Value = '  -6.91%  '
$ws.Range("E11").Value = '  -5.54%  '
$ws.Range("E12").Value = '  +0.94%  '
$ws.Range("E13").Value = '  -3.96%  '
$ws.Range("E14").Value = '  -6.22%  '
$ws.Range("E15").Value = '  -3.83%  '
$ws.Range("E16").Value = '  -5.65%  '
$ws.Range("E17").Value = '  -8.29%  '
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("E19").Value = '  -4.13%  '
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("E21").Value = '  -6.26%  '
$ws.Range("E22").Value = '  -4.77%  '
$ws.Range("E23").Value = '  -4.44%  '
$ws.Range("E24").Value = '  -1.35%  '
$ws.Range("E25").Value = '  -5.78%  '
$ws.Range("E26").Value = '  -5.67%  '
$ws.Range("E27").Value = '  -2.26%  '
$ws.Range("E28").Value = '  -5.54%  '
$ws.Range("E29").Value = '  -1.78%  '
$ws.Range("E30").Value = '  -3.95%  '
$ws.Range("E31").Value = '  -5.68%  '
$ws.Range("E32").Value = '  -1.77%  '
$ws.Range("E33").Value = '  -6.37%  '
$ws.Range("E34").Value = '  -5.60%  '
$ws.Range("E35").Value = '  -8.95%  '
$ws.Range("E36").Value = '  -9.39%  '
$ws.Range("E37").Value = '  +3.51%  '
$ws.Range("E38").Value = '  -4.87%  '
$ws.Range("E39").Value = '  +0.66%  '
$ws.Range("E40").Value = '  -4.63%  '
$ws.Range("E41").Value = '  -4.62%  '
$ws.Range("E42").Value = '  -8.54%  '
$ws.Range("E43").Value = '  -5.73%  '
$ws.Range("E44").Value = '  -6.07%  '
$ws.Range("E45").Value = '  -3.33%  '
$ws.Range("E46").Value = '  -3.84%  '
$ws.Range("E47").Value = '  -5.12%  '
$ws.Range("E48").Value = '  -4.47%  '
$ws.Range("E49").Value = '  -2.93%  '
$ws.Range("E50").Value = '  -5.86%  '
$ws.Range("E51").Value = '  -5.65%  '
